$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

$tbl.Cell(1,1).Range.Text = "17+4="
$tbl.Cell(1,2).Range.Text = "77-31="
$tbl.Cell(1,3).Range.Text = "59+0="
$tbl.Cell(1,4).Range.Text = "82-12="
$tbl.Cell(1,5).Range.Text = "98-93="
$tbl.Cell(2,1).Range.Text = "65-16="
$tbl.Cell(2,2).Range.Text = "60-32="
$tbl.Cell(2,3).Range.Text = "26+53="
$tbl.Cell(2,4).Range.Text = "75-28="
$tbl.Cell(2,5).Range.Text = "46+39="
$tbl.Cell(3,1).Range.Text = "44-2="
$tbl.Cell(3,2).Range.Text = "70-65="
$tbl.Cell(3,3).Range.Text = "43+47="
$tbl.Cell(3,4).Range.Text = "6+76="
$tbl.Cell(3,5).Range.Text = "51-43="
$tbl.Cell(4,1).Range.Text = "42+31="
$tbl.Cell(4,2).Range.Text = "96-92="
$tbl.Cell(4,3).Range.Text = "78-30="
$tbl.Cell(4,4).Range.Text = "17+70="
$tbl.Cell(4,5).Range.Text = "1+61="
$tbl.Cell(5,1).Range.Text = "5+17="
$tbl.Cell(5,2).Range.Text = "80-4="
$tbl.Cell(5,3).Range.Text = "10+68="
$tbl.Cell(5,4).Range.Text = "75-25="
$tbl.Cell(5,5).Range.Text = "99-50="
$tbl.Cell(6,1).Range.Text = "37+28="
$tbl.Cell(6,2).Range.Text = "82-81="
$tbl.Cell(6,3).Range.Text = "93-51="
$tbl.Cell(6,4).Range.Text = "52+32="
$tbl.Cell(6,5).Range.Text = "75+20="
$tbl.Cell(7,1).Range.Text = "51+38="
$tbl.Cell(7,2).Range.Text = "34+30="
$tbl.Cell(7,3).Range.Text = "97-39="
$tbl.Cell(7,4).Range.Text = "27+38="
$tbl.Cell(7,5).Range.Text = "57-0="
$tbl.Cell(8,1).Range.Text = "1+30="
$tbl.Cell(8,2).Range.Text = "67-7="
$tbl.Cell(8,3).Range.Text = "12+24="
$tbl.Cell(8,4).Range.Text = "9+83="
$tbl.Cell(8,5).Range.Text = "30+0="
$tbl.Cell(9,1).Range.Text = "96-63="
$tbl.Cell(9,2).Range.Text = "66-65="
$tbl.Cell(9,3).Range.Text = "93-54="
$tbl.Cell(9,4).Range.Text = "91-15="
$tbl.Cell(9,5).Range.Text = "96-21="
$tbl.Cell(10,1).Range.Text = "3+11="
$tbl.Cell(10,2).Range.Text = "76+11="
$tbl.Cell(10,3).Range.Text = "26-11="
$tbl.Cell(10,4).Range.Text = "64+13="
$tbl.Cell(10,5).Range.Text = "41+5="
$tbl.Cell(11,1).Range.Text = "42+31="
$tbl.Cell(11,2).Range.Text = "74-27="
$tbl.Cell(11,3).Range.Text = "29+0="
$tbl.Cell(11,4).Range.Text = "66-48="
$tbl.Cell(11,5).Range.Text = "11+10="
$tbl.Cell(12,1).Range.Text = "79-64="
$tbl.Cell(12,2).Range.Text = "34+54="
$tbl.Cell(12,3).Range.Text = "40-0="
$tbl.Cell(12,4).Range.Text = "57-11="
$tbl.Cell(12,5).Range.Text = "14-11="
$tbl.Cell(13,1).Range.Text = "31+2="
$tbl.Cell(13,2).Range.Text = "70-64="
$tbl.Cell(13,3).Range.Text = "61-21="
$tbl.Cell(13,4).Range.Text = "8+19="
$tbl.Cell(13,5).Range.Text = "46-35="
$tbl.Cell(14,1).Range.Text = "80-73="
$tbl.Cell(14,2).Range.Text = "5+69="
$tbl.Cell(14,3).Range.Text = "22+1="
$tbl.Cell(14,4).Range.Text = "94-64="
$tbl.Cell(14,5).Range.Text = "60+35="
$tbl.Cell(15,1).Range.Text = "9+19="
$tbl.Cell(15,2).Range.Text = "82-65="
$tbl.Cell(15,3).Range.Text = "95-83="
$tbl.Cell(15,4).Range.Text = "66-21="
$tbl.Cell(15,5).Range.Text = "11+72="
$tbl.Cell(16,1).Range.Text = "17+73="
$tbl.Cell(16,2).Range.Text = "12+56="
$tbl.Cell(16,3).Range.Text = "76+12="
$tbl.Cell(16,4).Range.Text = "40+12="
$tbl.Cell(16,5).Range.Text = "56-54="
$tbl.Cell(17,1).Range.Text = "26+8="
$tbl.Cell(17,2).Range.Text = "47+6="
$tbl.Cell(17,3).Range.Text = "58+2="
$tbl.Cell(17,4).Range.Text = "19-13="
$tbl.Cell(17,5).Range.Text = "26+49="
$tbl.Cell(18,1).Range.Text = "89-62="
$tbl.Cell(18,2).Range.Text = "2+61="
$tbl.Cell(18,3).Range.Text = "92-64="
$tbl.Cell(18,4).Range.Text = "93-56="
$tbl.Cell(18,5).Range.Text = "20+36="
$tbl.Cell(19,1).Range.Text = "20+13="
$tbl.Cell(19,2).Range.Text = "46-40="
$tbl.Cell(19,3).Range.Text = "27-18="
$tbl.Cell(19,4).Range.Text = "87-40="
$tbl.Cell(19,5).Range.Text = "56-44="
$tbl.Cell(20,1).Range.Text = "21-14="
$tbl.Cell(20,2).Range.Text = "64-17="
$tbl.Cell(20,3).Range.Text = "38-9="
$tbl.Cell(20,4).Range.Text = "27+70="
$tbl.Cell(20,5).Range.Text = "96-56="
